# Add new LeetCode exercises to the tracking sheet:
#  - row 28: k closest points to origin (array, heap)
#  - row 29: longest substring without repeating characters (string, array) + link
#  - row 30: 3sum (arrray, two pointers) + link
# Also widen column A to fit the new, longer problem names and update the
# active selection to the last entered cell.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 29 entered first (matches author's original typing order)
$ws.Cells.Item(29, 1).Value = "longest substring without repeating characters"
$ws.Cells.Item(29, 2).Value = "string"
$ws.Cells.Item(29, 3).Value = "array"
$ws.Cells.Item(29, 5).Value = "https://leetcode.com/problems/longest-substring-without-repeating-characters/"

# Row 28: k closest points to origin
$ws.Cells.Item(28, 1).Value = "k closest points to origin"
$ws.Cells.Item(28, 2).Value = "array"
$ws.Cells.Item(28, 3).Value = "heap"

# Row 30: 3sum
$ws.Cells.Item(30, 1).Value = "3sum"
$ws.Cells.Item(30, 5).Value = "https://leetcode.com/problems/3sum/"
$ws.Cells.Item(30, 2).Value = "arrray"
$ws.Cells.Item(30, 3).Value = "two pointers"

# Widen column A so the new, longer exercise names fit (best-fit autosize)
$ws.Columns.Item(1).AutoFit()

# Update view state: scroll so row 13 is at top, select last-entered cell D30
$ws.Range("A13").Select() | Out-Null
$ws.Range("D30").Select() | Out-Null
